# Weekly update: insert a new record row for "Coliflor" (Macroferia Regional
# de Talca) above the existing row 253, pushing the rest of the series down
# by one row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 253:265 down to 254:266 by inserting a new blank row at 253.
$ws.Rows("253:253").Insert()

# Populate the new row 253 with the latest weekly record, mirroring the
# format/layout of the surrounding rows.
$ws.Range("A253").Value = 5
$ws.Range("B253").Value = "Macroferia Regional de Talca"
$ws.Range("C253").Value = "Maule"
$ws.Range("D253").Value = 44753
$ws.Range("D253").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E253").Value = 7
$ws.Range("F253").Value = 100112008
$ws.Range("G253").Value = "Coliflor"
$ws.Range("H253").Value = "Sin especificar"
$ws.Range("I253").Value = "Primera"
$ws.Range("J253").Value = 3000
$ws.Range("K253").Value = 1000
$ws.Range("L253").Value = 1000
$ws.Range("M253").Value = 1000
$ws.Range("N253").Value = "`$/unidad"
$ws.Range("O253").Value = "Región del Maule"
$ws.Range("P253").Value = 1000
$ws.Range("Q253").Value = 1
$ws.Range("R253").Value = "Hortaliza"
